$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "SAN DIEGO AREA TOTALS" label moves from B2 to A2 (keeps its style).
$ws.Range("B2").Copy($ws.Range("A2"))

# B2 is cleared completely (contents + formatting) and re-populated with
# the "Totals" label (re-using the existing shared string), taking on the
# default (unstyled) cell format.
$ws.Range("B2").Clear()
$ws.Range("B2").Value = "Totals"

# Column A now holds the same kind of port-name text that column B already
# had, so it widens to (approximately) the same best-fit width B already
# has.
$ws.Columns("A:A").ColumnWidth = 21.83

# The active selection becomes the whole of column A.
$ws.Range("A1:A1048576").Select()

$wb.Save()
